# Replace end-to-end architecture: remove the "Step-by-Step Architecture"
# section (heading + "1. Data Collection Module" subsection) that sits
# between two blank paragraphs near the end of the document.

$d = $word.ActiveDocument

$paras = $d.Paragraphs
$sectionStart = $null
$sectionEnd = $null

for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($sectionStart -eq $null -and $t -match "Step-by-Step Architecture") {
        $sectionStart = $p
    }
    if ($sectionStart -ne $null -and $t -match "Works at scale with low overhead") {
        $sectionEnd = $p
        break
    }
}

if ($sectionStart -ne $null -and $sectionEnd -ne $null) {
    $r = $d.Range($sectionStart.Range.Start, $sectionEnd.Range.End)
    $r.Delete()
}
